$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 21, shifting existing rows 21-24 down to 22-25
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with the new record's data
$ws.Cells.Item(21, 1).Value = 10
$ws.Cells.Item(21, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(21, 3).Value = "La Araucanía"
$ws.Cells.Item(21, 4).Value = 44522
$ws.Cells.Item(21, 5).Value = 9
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100103
$ws.Cells.Item(21, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(21, 9).Value = 100103003
$ws.Cells.Item(21, 10).Value = "Damasco"
$ws.Cells.Item(21, 11).Value = "Castle Brite"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 30
$ws.Cells.Item(21, 14).Value = 35000
$ws.Cells.Item(21, 15).Value = 35000
$ws.Cells.Item(21, 16).Value = 35000
$ws.Cells.Item(21, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(21, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(21, 19).Value = 2333
$ws.Cells.Item(21, 20).Value = 15
